$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F2 changes from FAIL to PASS
$ws.Range("F2").Value = "PASS"

# F5 and F6 get new PASS values
$ws.Range("F5").Value = "PASS"
$ws.Range("F6").Value = "PASS"

# New rows 7-9
$ws.Range("C7").Value = "kuqnabgrrdhrdpd@gmail.com"
$ws.Range("D7").Value = "vhtvsMFBFZ5"
$ws.Range("E7").Value = "pass"
$ws.Range("F7").Value = "PASS"

$ws.Range("C8").Value = "uuvqmhzknkmjnpc@gmail.com"
$ws.Range("D8").Value = "vytszXNSCT5"
$ws.Range("E8").Value = "pass"
$ws.Range("F8").Value = "PASS"

$ws.Range("C9").Value = "ycwlrrnwcgpzsgf@gmail.com"
$ws.Range("D9").Value = "fejnwTVDEH5"
$ws.Range("E9").Value = "pass"
$ws.Range("F9").Value = "PASS"
